$wb = $excel.ActiveWorkbook

# --- Update "Ready for handoff" -> "Handback transform failed" (all usages) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"

# --- Error Detail column (P) gets new handback-mismatch messages ---
$wsZhCn.Range("P3").Value = "Handback file name: 0rxl2bg2.h5h is different with handoff file name: 501e4453-fef1-4f33-876e-1aa471f1a87f.9f98cda53b07940dd3e906ae422ce59e6664f5ea.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: 0rxl2bg2.h5h is different with handoff file name: 501e4453-fef1-4f33-876e-1aa471f1a87f.9f98cda53b07940dd3e906ae422ce59e6664f5ea.de-de."

# --- Widen Error Detail column (16 = P) from ~13.75 to 40 characters ---
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
